# Committing stabilized code for change material Nav as per current flow
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (Create_Material_Global_Local_NAV test case) updates
$ws.Range("C4").Value = "Y"
$ws.Range("E4").Value = "Test Create Automation"
$ws.Range("G4").Value = "CMG0472"
$ws.Range("K4").Value = "G"
$ws.Range("L4").Value = "YROH"
$ws.Range("N4").Value = "G"
$ws.Range("W4").Value = "G"
$ws.Range("X4").Value = "G"
$ws.Range("Y4").Value = "G"
$ws.Range("Z4").Value = "G"
$ws.Range("AE4").Value = "G"
$ws.Range("AF4").Value = "G"
$ws.Range("AG4").Value = "G"
$ws.Range("AN4").Value = "G"
$ws.Range("AO4").Value = "G"

# Row 5 (Create_Material_Global_and_Local_for_JDE test case) update
$ws.Range("C5").Value = "N"

# Leave the active selection where the edits were made
$ws.Range("D4").Select()
